$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove erroneous naive forecaster values in C2 and C3 (bug fix)
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Corrected floating point precision of forecast vectors
$ws.Range("C4").Value = -0.01587181126744275
$ws.Range("C5").Value = -0.02256889165886955
$ws.Range("E5").Value = 0.1850158025575199
$ws.Range("C6").Value = 0.09611428386595566
$ws.Range("E6").Value = -0.0461580488825164
$ws.Range("C8").Value = -0.001350220946472191
$ws.Range("E9").Value = -0.2383077634183106
$ws.Range("C10").Value = -0.5761528471665334
$ws.Range("E10").Value = -0.2318455351884796
$ws.Range("E11").Value = -0.2555440101933759
$ws.Range("C14").Value = -0.4278219446121501
$ws.Range("C15").Value = -1.026566979837429
$ws.Range("E15").Value = -1.192393303545602
$ws.Range("C17").Value = 0.4636049209196802
$ws.Range("C18").Value = 0.6216390921348403
$ws.Range("E18").Value = -0.1830321574487592
$ws.Range("C19").Value = -0.6768900623516871
